# -----------------------------------------------------------------------
# Rework Sheet1's test-case table (registration-form scenarios) and add
# five new sheets, each covering a "missing field" validation scenario.
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Sheet1: update the existing rows ----------------------------------
# Row 2 (User1) is already correct, no change needed.

# Row 3 (User2): password now matches confirmPassword -> success
$ws1.Range("D3").Value = "qwerty"
$ws1.Range("E3").ClearFormats()
$ws1.Range("E3").Value = "User has been successfully added."

# Row 4: was blank fullName "User5" row -> becomes "User3" / mismatch case
$ws1.Range("A4").Value = "User3"
$ws1.Range("D4").Value = "qwert"
$ws1.Range("E4").Value = "Passwords don't match."

# Row 5: was blank email "User6" row -> becomes "email4" / mismatch case
$ws1.Range("B5").Value = "email4"
$ws1.Range("D5").Value = "qwert"
$ws1.Range("E5").Value = "Passwords don't match."

# Remove the old rows 6 and 7 (User5 / User6 scenarios no longer needed)
$ws1.Rows("6:7").Delete()

# Column C got a touch wider
# (ColumnWidth goes through a 6px/char digit-width model on this host, so we
# solve for the "characters" input that reproduces the stored OOXML width
# as closely as possible.)
function ToColWidth([double]$target) {
    $p = [Math]::Round($target * 6)
    return ($p - 5) / 6
}
$ws1.Columns("C").ColumnWidth = ToColWidth 11.85546875

# Update the active selection/cursor to reflect the new last row
$ws1.Range("B6").Select()

# -------------------------------------------------------------------------
# Helper data used to build the five new "missing field" scenario sheets.
# Each sheet tests the registration form with one field left blank and
# shows the expected "Incomplete fields." result.
# -------------------------------------------------------------------------

function Add-ScenarioSheet {
    param(
        [string]$Name,
        [string[]]$Headers,
        [string[]]$Row2,
        [string[]]$Row3,
        [double[]]$ColWidths
    )

    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $last)
    $ws.Name = $Name

    for ($i = 0; $i -lt $Headers.Count; $i++) {
        $col = [char](65 + $i)
        $ws.Range("$col`1").Value = $Headers[$i]
    }
    for ($i = 0; $i -lt $Row2.Count; $i++) {
        $col = [char](65 + $i)
        $ws.Range("$col`2").Value = $Row2[$i]
    }
    for ($i = 0; $i -lt $Row3.Count; $i++) {
        $col = [char](65 + $i)
        $ws.Range("$col`3").Value = $Row3[$i]
    }

    for ($i = 0; $i -lt $ColWidths.Count; $i++) {
        $col = [char](65 + $i)
        $ws.Columns($col).ColumnWidth = ToColWidth $ColWidths[$i]
    }

    return $ws
}

# Sheet2: missing fullName
$s2 = Add-ScenarioSheet "Sheet2" `
    @("email","password","confirmPassword","Expected Value","Actual Value","Result") `
    @("email1","qwerty","qwerty","Incomplete fields.") `
    @("email2","qwerty","qwerty","Incomplete fields.") `
    @(16.42578125,31.7109375,12.140625,6.5703125)
$s2.Range("C1").Select()

# Sheet3: missing email
$s3 = Add-ScenarioSheet "Sheet3" `
    @("fullName","password","confirmPassword","Expected Value","Actual Value","Result") `
    @("User1","qwerty","qwerty","Incomplete fields.") `
    @("User2","qwerty","qwerty","Incomplete fields.") `
    @(9.28515625,9.42578125,16.42578125,17.42578125,12.140625,6.5703125)

# Sheet4: missing password
$s4 = Add-ScenarioSheet "Sheet4" `
    @("fullName","email","confirmPassword","Expected Value","Actual Value","Result") `
    @("User1","email1","qwerty","Incomplete fields.") `
    @("User2","email2","qwerty","Incomplete fields.") `
    @(9.28515625,9.28515625,16.42578125,17.42578125,12.140625,6.5703125)

# Sheet5: missing confirmPassword
$s5 = Add-ScenarioSheet "Sheet5" `
    @("fullName","email","password","Expected Value","Actual Value","Result") `
    @("User1","email1","qwerty","Incomplete fields.") `
    @("User2","email2","qwerty","Incomplete fields.") `
    @(9.28515625,7,9.42578125,17.42578125,12.140625,6.5703125)

# Sheet6: everything missing
$s6 = Add-ScenarioSheet "Sheet6" `
    @("Expected Value","Actual Value","Result") `
    @("Incomplete fields.") `
    @("Incomplete fields.") `
    @(17.42578125,12.140625)

# Re-select Sheet2 as the active tab, matching the new workbook view.
$s2.Activate()
